# Auto-generated edit script for cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'27.284.19"
$ws.Cells.Item(2, 5).Value = "  -0.61%  "

$ws.Cells.Item(3, 4).Value = "'1.708.75"
$ws.Cells.Item(3, 5).Value = "  -0.70%  "

$ws.Cells.Item(4, 4).Value = "'1.003"
$ws.Cells.Item(4, 5).Value = "  +0.08%  "

$ws.Cells.Item(5, 4).Value = "'223.72"
$ws.Cells.Item(5, 5).Value = "  -2.37%  "

$ws.Cells.Item(6, 4).Value = "'0.5285"
$ws.Cells.Item(6, 5).Value = "  -2.33%  "

$ws.Cells.Item(7, 4).Value = "'1.003"
$ws.Cells.Item(7, 5).Value = "  +0.10%  "

$ws.Cells.Item(8, 4).Value = "'0.2638"
$ws.Cells.Item(8, 5).Value = "  -4.46%  "

$ws.Cells.Item(9, 4).Value = "'0.06538"
$ws.Cells.Item(9, 5).Value = "  -2.99%  "

$ws.Cells.Item(10, 4).Value = "'20.91"
$ws.Cells.Item(10, 5).Value = "  -2.94%  "

$ws.Cells.Item(11, 4).Value = "'0.07630"
$ws.Cells.Item(11, 5).Value = "  -2.11%  "

$ws.Cells.Item(12, 4).Value = "'4.564"
$ws.Cells.Item(12, 5).Value = "  -3.30%  "

$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13, 4).Value = "'1.708.32"
$ws.Cells.Item(13, 5).Value = "  -1.03%  "

$ws.Cells.Item(14, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(14, 4).Value = "'1.945.10"
$ws.Cells.Item(14, 5).Value = "  -0.63%  "

$ws.Cells.Item(15, 5).Value = "  -3.92%  "

$ws.Cells.Item(16, 4).Value = "'0.0₅8175"
$ws.Cells.Item(16, 5).Value = "  -2.51%  "

$ws.Cells.Item(17, 4).Value = "'67.19"
$ws.Cells.Item(17, 5).Value = "  -2.15%  "

$ws.Cells.Item(18, 4).Value = "'27.272.37"
$ws.Cells.Item(18, 5).Value = "  -0.56%  "

$ws.Cells.Item(19, 4).Value = "'215.55"
$ws.Cells.Item(19, 5).Value = "  +2.54%  "

$ws.Cells.Item(20, 4).Value = "'1.003"
$ws.Cells.Item(20, 5).Value = "  +0.09%  "

$ws.Cells.Item(21, 4).Value = "'4.678"
$ws.Cells.Item(21, 5).Value = "  -2.75%  "

$ws.Cells.Item(22, 4).Value = "'10.45"
$ws.Cells.Item(22, 5).Value = "  -4.09%  "

$ws.Cells.Item(23, 4).Value = "'5.954"
$ws.Cells.Item(23, 5).Value = "  -4.58%  "

$ws.Cells.Item(24, 4).Value = "'1.004"
$ws.Cells.Item(24, 5).Value = "  +0.19%  "

$ws.Cells.Item(25, 4).Value = "'142.51"
$ws.Cells.Item(25, 5).Value = "  -2.79%  "

$ws.Cells.Item(26, 4).Value = "'1.753"
$ws.Cells.Item(26, 5).Value = "  +7.36%  "

$ws.Cells.Item(27, 4).Value = "'0.1216"
$ws.Cells.Item(27, 5).Value = "  -2.90%  "

$ws.Cells.Item(28, 4).Value = "'7.261"
$ws.Cells.Item(28, 5).Value = "  -2.42%  "

$ws.Cells.Item(29, 4).Value = "'16.29"
$ws.Cells.Item(29, 5).Value = "  -3.51%  "

$ws.Cells.Item(30, 4).Value = "'0.05374"
$ws.Cells.Item(30, 5).Value = "  -3.87%  "

$ws.Cells.Item(31, 5).Value = "  -1.48%  "

$ws.Cells.Item(32, 4).Value = "'3.483"
$ws.Cells.Item(32, 5).Value = "  -4.51%  "

$ws.Cells.Item(33, 4).Value = "'3.411"
$ws.Cells.Item(33, 5).Value = "  -2.84%  "

$ws.Cells.Item(34, 5).Value = "  +0.16%  "

$ws.Cells.Item(35, 4).Value = "'2.869"
$ws.Cells.Item(35, 5).Value = "  +0.75%  "

$ws.Cells.Item(36, 2).Value = "ARBITRUM"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(36, 4).Value = "'0.9510"
$ws.Cells.Item(36, 5).Value = "  -2.54%  "

$ws.Cells.Item(37, 2).Value = "HuobiToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(37, 4).Value = "'2.422"
$ws.Cells.Item(37, 5).Value = "  -0.79%  "

$ws.Cells.Item(38, 4).Value = "'0.5860"
$ws.Cells.Item(38, 5).Value = "  +0.43%  "

$ws.Cells.Item(39, 5).Value = "  -0.76%  "

$ws.Cells.Item(40, 4).Value = "'5.861"
$ws.Cells.Item(40, 5).Value = "  +0.03%  "

$ws.Cells.Item(41, 4).Value = "'1.003"
$ws.Cells.Item(41, 5).Value = "  +0.14%  "

$ws.Cells.Item(42, 2).Value = "Maker"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(42, 4).Value = "'1.039.41"
$ws.Cells.Item(42, 5).Value = "  -0.23%  "

$ws.Cells.Item(43, 2).Value = "TrustWalletToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(43, 4).Value = "'0.8393"
$ws.Cells.Item(43, 5).Value = "  -0.05%  "

$ws.Cells.Item(44, 4).Value = "'101.00"
$ws.Cells.Item(44, 5).Value = "  -1.45%  "

$ws.Cells.Item(45, 4).Value = "'1.852.86"
$ws.Cells.Item(45, 5).Value = "  -0.49%  "

$ws.Cells.Item(46, 5).Value = "  +5.66%  "

$ws.Cells.Item(47, 4).Value = "'57.96"
$ws.Cells.Item(47, 5).Value = "  -2.90%  "

$ws.Cells.Item(48, 4).Value = "'0.4489"
$ws.Cells.Item(48, 5).Value = "  +2.16%  "

$ws.Cells.Item(49, 5).Value = "  +0.28%  "

$ws.Cells.Item(50, 4).Value = "'8.076"
$ws.Cells.Item(50, 5).Value = "  -1.32%  "

$ws.Cells.Item(51, 4).Value = "'0.05239"
$ws.Cells.Item(51, 5).Value = "  -0.73%  "

